$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: copy style from E1 (bold header style) then set the text
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "time_taken"

# time_taken values for rows 2..280
$timeTaken = @(
    "2021-10-05 10:51:29.198897",
    "2021-10-05 10:51:29.198909",
    "2021-10-05 10:51:29.198912",
    "2021-10-05 10:51:29.198915",
    "2021-10-05 10:51:29.198918",
    "2021-10-05 10:51:29.198921",
    "2021-10-05 10:51:29.198924",
    "2021-10-05 10:51:29.198927",
    "2021-10-05 10:51:29.198930",
    "2021-10-05 10:51:29.198932",
    "2021-10-05 10:51:29.198935",
    "2021-10-05 10:51:29.198938",
    "2021-10-05 10:51:29.198940",
    "2021-10-05 10:51:29.198943",
    "2021-10-05 10:51:29.198946",
    "2021-10-05 10:51:29.198948",
    "2021-10-05 10:51:29.198951",
    "2021-10-05 10:51:29.198954",
    "2021-10-05 10:51:29.198957",
    "2021-10-05 10:51:29.198959",
    "2021-10-05 10:51:29.198962",
    "2021-10-05 10:51:29.198965",
    "2021-10-05 10:51:29.198968",
    "2021-10-05 10:51:29.198970",
    "2021-10-05 10:51:29.198973",
    "2021-10-05 10:51:29.198976",
    "2021-10-05 10:51:29.198979",
    "2021-10-05 10:51:29.198982",
    "2021-10-05 10:51:29.198985",
    "2021-10-05 10:51:29.198988",
    "2021-10-05 10:51:29.198991",
    "2021-10-05 10:51:29.198993",
    "2021-10-05 10:51:29.198996",
    "2021-10-05 10:51:29.198999",
    "2021-10-05 10:51:29.199002",
    "2021-10-05 10:51:29.199005",
    "2021-10-05 10:51:29.199007",
    "2021-10-05 10:51:29.199010",
    "2021-10-05 10:51:29.199013",
    "2021-10-05 10:51:29.199016",
    "2021-10-05 10:51:29.199019",
    "2021-10-05 10:51:29.199022",
    "2021-10-05 10:51:29.199024",
    "2021-10-05 10:51:29.199027",
    "2021-10-05 10:51:29.199030",
    "2021-10-05 10:51:29.199033",
    "2021-10-05 10:51:29.199035",
    "2021-10-05 10:51:29.199038",
    "2021-10-05 10:51:29.199041",
    "2021-10-05 10:51:29.199043",
    "2021-10-05 10:51:29.199046",
    "2021-10-05 10:51:29.199049",
    "2021-10-05 10:51:29.199052",
    "2021-10-05 10:51:29.199055",
    "2021-10-05 10:51:29.199058",
    "2021-10-05 10:51:29.199060",
    "2021-10-05 10:51:29.199063",
    "2021-10-05 10:51:29.199065",
    "2021-10-05 10:51:29.199068",
    "2021-10-05 10:51:29.199071",
    "2021-10-05 10:51:29.199073",
    "2021-10-05 10:51:29.199076",
    "2021-10-05 10:51:29.199079",
    "2021-10-05 10:51:29.199081",
    "2021-10-05 10:51:29.199086",
    "2021-10-05 10:51:29.199089",
    "2021-10-05 10:51:29.199091",
    "2021-10-05 10:51:29.199094",
    "2021-10-05 10:51:29.199097",
    "2021-10-05 10:51:29.199100",
    "2021-10-05 10:51:29.199102",
    "2021-10-05 10:51:29.199105",
    "2021-10-05 10:51:29.199108",
    "2021-10-05 10:51:29.199111",
    "2021-10-05 10:51:29.199113",
    "2021-10-05 10:51:29.199116",
    "2021-10-05 10:51:29.199121",
    "2021-10-05 10:51:29.199124",
    "2021-10-05 10:51:29.199127",
    "2021-10-05 10:51:29.199129",
    "2021-10-05 10:51:29.199132",
    "2021-10-05 10:51:29.199135",
    "2021-10-05 10:51:29.199138",
    "2021-10-05 10:51:29.199140",
    "2021-10-05 10:51:29.199143",
    "2021-10-05 10:51:29.199146",
    "2021-10-05 10:51:29.199149",
    "2021-10-05 10:51:29.199151",
    "2021-10-05 10:51:29.199154",
    "2021-10-05 10:51:29.199157",
    "2021-10-05 10:51:29.199160",
    "2021-10-05 10:51:29.199162",
    "2021-10-05 10:51:29.199166",
    "2021-10-05 10:51:29.199169",
    "2021-10-05 10:51:29.199172",
    "2021-10-05 10:51:29.199175",
    "2021-10-05 10:51:29.199177",
    "2021-10-05 10:51:29.199180",
    "2021-10-05 10:51:29.199183",
    "2021-10-05 10:51:29.199185",
    "2021-10-05 10:51:29.199188",
    "2021-10-05 10:51:29.199191",
    "2021-10-05 10:51:29.199193",
    "2021-10-05 10:51:29.199196",
    "2021-10-05 10:51:29.199199",
    "2021-10-05 10:51:29.199202",
    "2021-10-05 10:51:29.199205",
    "2021-10-05 10:51:29.199207",
    "2021-10-05 10:51:29.199212",
    "2021-10-05 10:51:29.199215",
    "2021-10-05 10:51:29.199219",
    "2021-10-05 10:51:29.199222",
    "2021-10-05 10:51:29.199224",
    "2021-10-05 10:51:29.199227",
    "2021-10-05 10:51:29.199230",
    "2021-10-05 10:51:29.199232",
    "2021-10-05 10:51:29.199235",
    "2021-10-05 10:51:29.199237",
    "2021-10-05 10:51:29.199240",
    "2021-10-05 10:51:29.199243",
    "2021-10-05 10:51:29.199245",
    "2021-10-05 10:51:29.199248",
    "2021-10-05 10:51:29.199251",
    "2021-10-05 10:51:29.199253",
    "2021-10-05 10:51:29.199256",
    "2021-10-05 10:51:29.199259",
    "2021-10-05 10:51:29.199261",
    "2021-10-05 10:51:29.199264",
    "2021-10-05 10:51:29.199268",
    "2021-10-05 10:51:29.199272",
    "2021-10-05 10:51:29.199274",
    "2021-10-05 10:51:29.199277",
    "2021-10-05 10:51:29.199280",
    "2021-10-05 10:51:29.199283",
    "2021-10-05 10:51:29.199285",
    "2021-10-05 10:51:29.199288",
    "2021-10-05 10:51:29.199291",
    "2021-10-05 10:51:29.199294",
    "2021-10-05 10:51:29.199297",
    "2021-10-05 10:51:29.199299",
    "2021-10-05 10:51:29.199302",
    "2021-10-05 10:51:29.199305",
    "2021-10-05 10:51:29.199307",
    "2021-10-05 10:51:29.199310",
    "2021-10-05 10:51:29.199313",
    "2021-10-05 10:51:29.199315",
    "2021-10-05 10:51:29.199318",
    "2021-10-05 10:51:29.199321",
    "2021-10-05 10:51:29.199324",
    "2021-10-05 10:51:29.199327",
    "2021-10-05 10:51:29.199329",
    "2021-10-05 10:51:29.199332",
    "2021-10-05 10:51:29.199335",
    "2021-10-05 10:51:29.199337",
    "2021-10-05 10:51:29.199340",
    "2021-10-05 10:51:29.199343",
    "2021-10-05 10:51:29.199346",
    "2021-10-05 10:51:29.199348",
    "2021-10-05 10:51:29.199351",
    "2021-10-05 10:51:29.199354",
    "2021-10-05 10:51:29.199356",
    "2021-10-05 10:51:29.199359",
    "2021-10-05 10:51:29.199362",
    "2021-10-05 10:51:29.199365",
    "2021-10-05 10:51:29.199367",
    "2021-10-05 10:51:29.199370",
    "2021-10-05 10:51:29.199372",
    "2021-10-05 10:51:29.199375",
    "2021-10-05 10:51:29.199378",
    "2021-10-05 10:51:29.199380",
    "2021-10-05 10:51:29.199383",
    "2021-10-05 10:51:29.199386",
    "2021-10-05 10:51:29.199390",
    "2021-10-05 10:51:29.199393",
    "2021-10-05 10:51:29.199396",
    "2021-10-05 10:51:29.199399",
    "2021-10-05 10:51:29.199402",
    "2021-10-05 10:51:29.199404",
    "2021-10-05 10:51:29.199407",
    "2021-10-05 10:51:29.199410",
    "2021-10-05 10:51:29.199412",
    "2021-10-05 10:51:29.199415",
    "2021-10-05 10:51:29.199418",
    "2021-10-05 10:51:29.199420",
    "2021-10-05 10:51:29.199423",
    "2021-10-05 10:51:29.199426",
    "2021-10-05 10:51:29.199428",
    "2021-10-05 10:51:29.199431",
    "2021-10-05 10:51:29.199434",
    "2021-10-05 10:51:29.199436",
    "2021-10-05 10:51:29.199439",
    "2021-10-05 10:51:29.199442",
    "2021-10-05 10:51:29.199444",
    "2021-10-05 10:51:29.199447",
    "2021-10-05 10:51:29.199450",
    "2021-10-05 10:51:29.199452",
    "2021-10-05 10:51:29.199455",
    "2021-10-05 10:51:29.199458",
    "2021-10-05 10:51:29.199460",
    "2021-10-05 10:51:29.199463",
    "2021-10-05 10:51:29.199466",
    "2021-10-05 10:51:29.199469",
    "2021-10-05 10:51:29.199472",
    "2021-10-05 10:51:29.199474",
    "2021-10-05 10:51:29.199477",
    "2021-10-05 10:51:29.199480",
    "2021-10-05 10:51:29.199482",
    "2021-10-05 10:51:29.199485",
    "2021-10-05 10:51:29.199487",
    "2021-10-05 10:51:29.199490",
    "2021-10-05 10:51:29.199493",
    "2021-10-05 10:51:29.199495",
    "2021-10-05 10:51:29.199498",
    "2021-10-05 10:51:29.199501",
    "2021-10-05 10:51:29.199503",
    "2021-10-05 10:51:29.199506",
    "2021-10-05 10:51:29.199509",
    "2021-10-05 10:51:29.199511",
    "2021-10-05 10:51:29.199514",
    "2021-10-05 10:51:29.199517",
    "2021-10-05 10:51:29.199520",
    "2021-10-05 10:51:29.199545",
    "2021-10-05 10:51:29.199554",
    "2021-10-05 10:51:29.199559",
    "2021-10-05 10:51:29.199563",
    "2021-10-05 10:51:29.199568",
    "2021-10-05 10:51:29.199573",
    "2021-10-05 10:51:29.199576",
    "2021-10-05 10:51:29.199580",
    "2021-10-05 10:51:29.199583",
    "2021-10-05 10:51:29.199586",
    "2021-10-05 10:51:29.199589",
    "2021-10-05 10:51:29.199598",
    "2021-10-05 10:51:29.199603",
    "2021-10-05 10:51:29.199606",
    "2021-10-05 10:51:29.199609",
    "2021-10-05 10:51:29.199612",
    "2021-10-05 10:51:29.199614",
    "2021-10-05 10:51:29.199617",
    "2021-10-05 10:51:29.199620",
    "2021-10-05 10:51:29.199625",
    "2021-10-05 10:51:29.199629",
    "2021-10-05 10:51:29.199633",
    "2021-10-05 10:51:29.199637",
    "2021-10-05 10:51:29.199639",
    "2021-10-05 10:51:29.199643",
    "2021-10-05 10:51:29.199646",
    "2021-10-05 10:51:29.199649",
    "2021-10-05 10:51:29.199652",
    "2021-10-05 10:51:29.199654",
    "2021-10-05 10:51:29.199657",
    "2021-10-05 10:51:29.199659",
    "2021-10-05 10:51:29.199662",
    "2021-10-05 10:51:29.199664",
    "2021-10-05 10:51:29.199667",
    "2021-10-05 10:51:29.199670",
    "2021-10-05 10:51:29.199672",
    "2021-10-05 10:51:29.199675",
    "2021-10-05 10:51:29.199677",
    "2021-10-05 10:51:29.199680",
    "2021-10-05 10:51:29.199683",
    "2021-10-05 10:51:29.199685",
    "2021-10-05 10:51:29.199688",
    "2021-10-05 10:51:29.199691",
    "2021-10-05 10:51:29.199693",
    "2021-10-05 10:51:29.199696",
    "2021-10-05 10:51:29.199698",
    "2021-10-05 10:51:29.199701",
    "2021-10-05 10:51:29.199705",
    "2021-10-05 10:51:29.199710",
    "2021-10-05 10:51:29.199715",
    "2021-10-05 10:51:29.199720",
    "2021-10-05 10:51:29.199736",
    "2021-10-05 10:51:29.199740",
    "2021-10-05 10:51:29.199743",
    "2021-10-05 10:51:29.199746",
    "2021-10-05 10:51:29.199749",
    "2021-10-05 10:51:29.199752",
    "2021-10-05 10:51:29.199754"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

$excel.CutCopyMode = 0
